$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight the P6:S6 (+S7) and U6 (+U7) header cells with new fill colours ---
# These are new distinct fills, which is why new cellXfs entries get created.
# OLE colours are R + G*256 + B*65536.
$purple = 213 + (213*256) + (255*65536)   # -> FFD5D5FF
$pink   = 255 + (221*256) + (255*65536)   # -> FFFFDDFF
$ws.Range("P6:S6").Interior.Color = $purple
$ws.Range("S7").Interior.Color = $purple
$ws.Range("U6").Interior.Color = $pink
$ws.Range("U7").Interior.Color = $pink

# --- Fix row 12 (Charles Martin): C12 should say "CP" like the other rows, ---
# --- with the same border/centred style as the other C-column cells. ---
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = "CP"

# --- Rename "Gallet Benjamin" -> "G Benjamin" (row 15) ---
$ws.Range("B15").Value = "G Benjamin"

# --- Add the new "taux horaire" test column E, copying D's look/style ---
$ws.Range("D8").Copy() | Out-Null
$ws.Range("E8:E15").PasteSpecial(-4122) | Out-Null

$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 170

$excel.CutCopyMode = 0
